$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells I16, I17, J17, I18, J18 hold numeric-looking text (coordinate lists /
# confidence scores) that must stay as TEXT, not be auto-coerced to numbers
# by COM's smart typing. Force the Text format before writing them.
$textCells = @("I16", "I17", "J17", "I18", "J18")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 16: image re-capture + slightly refined bounding box
$ws.Range("D16").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("I16").Value = "642,530,686,576"

# Row 17: image re-capture + refined bounding box + confidence bump
$ws.Range("D17").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("I17").Value = "794,481,831,526"
$ws.Range("J17").Value = "0.71"

# Row 18: different source image + refined bounding box + confidence change
$ws.Range("D18").Value = "image_20250808100711_ppp0.jpg"
$ws.Range("I18").Value = "1182,409,1232,451"
$ws.Range("J18").Value = "0.75"
